$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The DA_CE column (G) values for all data rows (2-65) are cleared out,
# leaving the cells blank while keeping the header in G1 ("DA_CE") intact.
$ws.Range("G2:G65").ClearContents()
